$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 316, shifting existing rows 316:333 down to 317:334.
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row 316 with the new weekly record.
$ws.Cells.Item(316, 1).Value = 3
$ws.Cells.Item(316, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(316, 3).Value = "Coquimbo"
$ws.Cells.Item(316, 4).Value = 44516
$ws.Cells.Item(316, 5).Value = 5
$ws.Cells.Item(316, 6).Value = 100112037
$ws.Cells.Item(316, 7).Value = "Cebollín"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 160
$ws.Cells.Item(316, 11).Value = 3000
$ws.Cells.Item(316, 12).Value = 3000
$ws.Cells.Item(316, 13).Value = 3000
$ws.Cells.Item(316, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(316, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(316, 16).Value = 83
$ws.Cells.Item(316, 17).Value = 36
$ws.Cells.Item(316, 18).Value = "Hortaliza"
